# Update the team-specific transition matrix (Winthrop_A) cell values
# to reflect additional simulated games (recomputed transition
# probabilities) per the commit: "added more games, sped up simulate
# game logic, and drafted optimization logic".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2025
$ws.Cells.Item(2, 3).Value = 0.5575
$ws.Cells.Item(2, 10).Value = 0.015
$ws.Cells.Item(2, 16).Value = 0.125
$ws.Cells.Item(2, 19).Value = 0.1
$ws.Cells.Item(3, 2).Value = 0.008888888888888889
$ws.Cells.Item(3, 3).Value = 0.004444444444444444
$ws.Cells.Item(3, 10).Value = 0.008888888888888889
$ws.Cells.Item(3, 16).Value = 0.7555555555555555
$ws.Cells.Item(3, 19).Value = 0.2222222222222222
$ws.Cells.Item(4, 10).Value = 0.07142857142857142
$ws.Cells.Item(4, 15).Value = 0.01428571428571429
$ws.Cells.Item(4, 16).Value = 0.6714285714285714
$ws.Cells.Item(4, 19).Value = 0.2428571428571429
$ws.Cells.Item(6, 2).Value = 0.07079646017699115
$ws.Cells.Item(6, 4).Value = 0.02212389380530973
$ws.Cells.Item(6, 5).Value = 0.004424778761061947
$ws.Cells.Item(6, 6).Value = 0.06637168141592921
$ws.Cells.Item(6, 10).Value = 0.4026548672566372
$ws.Cells.Item(6, 15).Value = 0.008849557522123894
$ws.Cells.Item(6, 17).Value = 0.1415929203539823
$ws.Cells.Item(6, 18).Value = 0.03982300884955752
$ws.Cells.Item(6, 19).Value = 0.2433628318584071
$ws.Cells.Item(7, 2).Value = 0.1311475409836066
$ws.Cells.Item(7, 4).Value = 0.0273224043715847
$ws.Cells.Item(7, 6).Value = 0.03278688524590164
$ws.Cells.Item(7, 10).Value = 0.1256830601092896
$ws.Cells.Item(7, 15).Value = 0.03278688524590164
$ws.Cells.Item(7, 17).Value = 0.1256830601092896
$ws.Cells.Item(7, 18).Value = 0.1092896174863388
$ws.Cells.Item(7, 19).Value = 0.4153005464480874
$ws.Cells.Item(8, 2).Value = 0.1614457831325301
$ws.Cells.Item(8, 4).Value = 0.02650602409638554
$ws.Cells.Item(8, 6).Value = 0.05301204819277108
$ws.Cells.Item(8, 10).Value = 0.108433734939759
$ws.Cells.Item(8, 15).Value = 0.02409638554216868
$ws.Cells.Item(8, 17).Value = 0.1566265060240964
$ws.Cells.Item(8, 18).Value = 0.09397590361445783
$ws.Cells.Item(8, 19).Value = 0.3759036144578313
$ws.Cells.Item(9, 2).Value = 0.1157894736842105
$ws.Cells.Item(9, 4).Value = 0.02105263157894737
$ws.Cells.Item(9, 5).Value = 0.005263157894736842
$ws.Cells.Item(9, 6).Value = 0.02631578947368421
$ws.Cells.Item(9, 10).Value = 0.1368421052631579
$ws.Cells.Item(9, 15).Value = 0.02105263157894737
$ws.Cells.Item(9, 17).Value = 0.2052631578947368
$ws.Cells.Item(9, 18).Value = 0.07368421052631578
$ws.Cells.Item(9, 19).Value = 0.3947368421052632
$ws.Cells.Item(10, 2).Value = 0.1423047177107502
$ws.Cells.Item(10, 4).Value = 0.03480278422273782
$ws.Cells.Item(10, 5).Value = 0.002320185614849188
$ws.Cells.Item(10, 6).Value = 0.05413766434648105
$ws.Cells.Item(10, 10).Value = 0.119876256767208
$ws.Cells.Item(10, 15).Value = 0.01160092807424594
$ws.Cells.Item(10, 17).Value = 0.1979891724671307
$ws.Cells.Item(10, 18).Value = 0.08816705336426914
$ws.Cells.Item(10, 19).Value = 0.3488012374323279
$ws.Cells.Item(11, 7).Value = 0.1134751773049645
$ws.Cells.Item(11, 10).Value = 0.07446808510638298
$ws.Cells.Item(11, 11).Value = 0.1773049645390071
$ws.Cells.Item(11, 12).Value = 0.6134751773049646
$ws.Cells.Item(11, 19).Value = 0.02127659574468085
$ws.Cells.Item(12, 7).Value = 0.7727272727272727
$ws.Cells.Item(12, 10).Value = 0.1590909090909091
$ws.Cells.Item(12, 11).Value = 0.01136363636363636
$ws.Cells.Item(12, 12).Value = 0.01704545454545454
$ws.Cells.Item(12, 19).Value = 0.03977272727272727
$ws.Cells.Item(13, 7).Value = 0.6216216216216216
$ws.Cells.Item(13, 10).Value = 0.3243243243243243
$ws.Cells.Item(13, 19).Value = 0.05405405405405406
$ws.Cells.Item(15, 6).Value = 0.03347280334728033
$ws.Cells.Item(15, 8).Value = 0.09623430962343096
$ws.Cells.Item(15, 9).Value = 0.07531380753138076
$ws.Cells.Item(15, 10).Value = 0.401673640167364
$ws.Cells.Item(15, 11).Value = 0.08368200836820083
$ws.Cells.Item(15, 13).Value = 0.004184100418410041
$ws.Cells.Item(15, 15).Value = 0.07112970711297072
$ws.Cells.Item(15, 19).Value = 0.2343096234309623
$ws.Cells.Item(16, 6).Value = 0.02264150943396226
$ws.Cells.Item(16, 8).Value = 0.1735849056603773
$ws.Cells.Item(16, 9).Value = 0.09056603773584905
$ws.Cells.Item(16, 10).Value = 0.3547169811320754
$ws.Cells.Item(16, 11).Value = 0.1056603773584906
$ws.Cells.Item(16, 13).Value = 0.01886792452830189
$ws.Cells.Item(16, 14).Value = 0.003773584905660377
$ws.Cells.Item(16, 15).Value = 0.0830188679245283
$ws.Cells.Item(16, 19).Value = 0.1471698113207547
$ws.Cells.Item(17, 6).Value = 0.05528846153846154
$ws.Cells.Item(17, 8).Value = 0.1682692307692308
$ws.Cells.Item(17, 9).Value = 0.0673076923076923
$ws.Cells.Item(17, 10).Value = 0.40625
$ws.Cells.Item(17, 11).Value = 0.07932692307692307
$ws.Cells.Item(17, 13).Value = 0.01682692307692308
$ws.Cells.Item(17, 14).Value = 0.002403846153846154
$ws.Cells.Item(17, 15).Value = 0.08653846153846154
$ws.Cells.Item(17, 19).Value = 0.1177884615384615
$ws.Cells.Item(18, 6).Value = 0.04166666666666666
$ws.Cells.Item(18, 8).Value = 0.1197916666666667
$ws.Cells.Item(18, 9).Value = 0.140625
$ws.Cells.Item(18, 10).Value = 0.4375
$ws.Cells.Item(18, 11).Value = 0.08854166666666667
$ws.Cells.Item(18, 13).Value = 0.02083333333333333
$ws.Cells.Item(18, 15).Value = 0.04166666666666666
$ws.Cells.Item(18, 19).Value = 0.109375
$ws.Cells.Item(19, 6).Value = 0.02834008097165992
$ws.Cells.Item(19, 8).Value = 0.2032388663967611
$ws.Cells.Item(19, 9).Value = 0.07692307692307693
$ws.Cells.Item(19, 10).Value = 0.3643724696356275
$ws.Cells.Item(19, 11).Value = 0.108502024291498
$ws.Cells.Item(19, 13).Value = 0.01781376518218623
$ws.Cells.Item(19, 15).Value = 0.07611336032388664
$ws.Cells.Item(19, 19).Value = 0.1246963562753036